# The "Presidencial / PASO" placeholder row is replaced with a full, real
# election record (GenCor2019 - Gobernador y Vice Gobernador de Córdoba).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "GenCor2019"
$ws.Range("B2").Value = "Elecciones Generales de la Provincia de Córdoba"
$ws.Range("C2").Value = "Gobernador y Vice Gobernador"

# D2/E2 hold an ISO date (2019-05-12) as plain text, not a date serial —
# entering it with a leading apostrophe keeps Excel from auto-converting it
# to a date value, and ClearFormats() drops the resulting quote-prefix style
# so the cell is left with no explicit style (matching the source data).
$ws.Range("D2").Value = "'2019-05-12"
$ws.Range("D2").ClearFormats()

$ws.Range("E2").Value = "'2019-05-12"
$ws.Range("E2").ClearFormats()
